$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 is a new batch row (batch 1345), essentially a duplicate of row 6
# (batch 1344) with the sort/second-sort/transaction numbers bumped to the
# new numeric batch id 1345 instead of the old "1344" text values.
$ws.Cells.Item(7,1).Value = 1345
$ws.Cells.Item(7,2).Value = 1345
$ws.Cells.Item(7,3).Value = "Batch No: 1344"
$ws.Cells.Item(7,4).Value = 1
$ws.Cells.Item(7,5).Value = 1345
$ws.Cells.Item(7,6).Value = 350
$ws.Cells.Item(7,7).Value = "P76       "
$ws.Cells.Item(7,8).Value = "2020 RAM 2500"

$ws.Cells.Item(7,9).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(7,9).Value = 45613

# These are numeric-looking strings that must stay stored as text (no
# explicit cell style though), so force text format, assign, then drop
# back to the default "Normal" style to avoid leaving a number format
# override behind.
$ws.Cells.Item(7,10).NumberFormat = "@"
$ws.Cells.Item(7,10).Value = "    231675"
$ws.Cells.Item(7,10).Style = "Normal"

$ws.Cells.Item(7,11).NumberFormat = "@"
$ws.Cells.Item(7,11).Value = "         1"
$ws.Cells.Item(7,11).Style = "Normal"

$ws.Cells.Item(7,12).NumberFormat = "@"
$ws.Cells.Item(7,12).Value = "      1001"

$ws.Cells.Item(7,13).NumberFormat = "@"
$ws.Cells.Item(7,13).Value = "    4"
$ws.Cells.Item(7,13).Style = "Normal"

$ws.Cells.Item(7,14).Value = 5201
$ws.Cells.Item(7,15).Value = 5201
$ws.Cells.Item(7,16).Value = 6036
$ws.Cells.Item(7,17).Value = 6036
$ws.Cells.Item(7,18).Value = 0.01
$ws.Cells.Item(7,19).Value = "R"
$ws.Cells.Item(7,20).Value = 5
$ws.Cells.Item(7,21).Value = 0.05
$ws.Cells.Item(7,22).Value = "PAM"

$ws.Cells.Item(7,33).Value = 25856.55
$ws.Cells.Item(7,36).Value = 1582913.2
$ws.Cells.Item(7,37).Value = "Total for Batch 1344:"
$ws.Cells.Item(7,38).Value = "EQ "
$ws.Cells.Item(7,39).Value = 0
$ws.Cells.Item(7,40).Value = 0
$ws.Cells.Item(7,41).Value = 1
$ws.Cells.Item(7,44).Value = 1
$ws.Cells.Item(7,45).Value = "2020 RAM 2500"
$ws.Cells.Item(7,46).Value = 0
$ws.Cells.Item(7,49).Value = "M"
$ws.Cells.Item(7,52).Value = "M"

# Moving to a new row shifts the active selection down to the next entry cell.
$null = $ws.Range("B8").Select()
